# Rbp4-Stra6.xlsx: refresh with new TPM numbers.
# - Row 2 (target cluster ECs) is recalculated and now targets FAPs.
# - Row 3 (target cluster FAPs) is recalculated and now targets MuSCs.
# - Row 4 (target cluster MuSCs) is dropped entirely (its data merged/superseded).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: target cluster ECs -> FAPs, with refreshed TPM-derived stats ---
$ws.Cells.Item(2, 4).Value = "FAPs"
$ws.Cells.Item(2, 11).Value = 3
$ws.Cells.Item(2, 12).Value = 1
$ws.Cells.Item(2, 13).Value = 1.003012
$ws.Cells.Item(2, 14).Value = 3.009036
$ws.Cells.Item(2, 15).Value = 0.419577131549034
$ws.Cells.Item(2, 16).Value = 0.419577131549034
$ws.Cells.Item(2, 17).Value = 9.356418571526667
$ws.Cells.Item(2, 18).Value = 84.20776714374001
$ws.Cells.Item(2, 19).Value = 0.419577131549034
$ws.Cells.Item(2, 20).Value = 0.419577131549034

# --- Row 3: target cluster FAPs -> MuSCs, with refreshed TPM-derived stats ---
$ws.Cells.Item(3, 4).Value = "MuSCs"
$ws.Cells.Item(3, 13).Value = 1.387518666666667
$ws.Cells.Item(3, 14).Value = 4.162555999999999
$ws.Cells.Item(3, 15).Value = 0.5804228684509659
$ws.Cells.Item(3, 16).Value = 0.5804228684509659
$ws.Cells.Item(3, 17).Value = 12.94322044117111
$ws.Cells.Item(3, 18).Value = 116.48898397054
$ws.Cells.Item(3, 19).Value = 0.5804228684509659
$ws.Cells.Item(3, 20).Value = 0.5804228684509659

# --- Row 4 (old target cluster MuSCs) no longer exists in the refreshed output ---
$ws.Rows(4).Delete()
